$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D against auto-numeric conversion by temporarily setting text format
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.374.41"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "1.877.18"

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "0.7138"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("D6").Value = "241.97"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "0.3117"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").Value = "0.07722"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("D10").Value = "25.19"
$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("D11").Value = "0.08383"
$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").Value = "1.891.37"
$ws.Range("E12").Value = "  +0.80%  "

$ws.Range("D13").Value = "5.241"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "0.7152"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").Value = "91.84"
$ws.Range("E15").Value = "  +1.11%  "

$ws.Range("D16").Value = "29.379.28"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").Value = "0.000008330"
$ws.Range("E17").Value = "  +6.29%  "

$ws.Range("D18").Value = "5.970"
$ws.Range("E18").Value = "  +1.76%  "

$ws.Range("D19").Value = "243.50"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "2.135.76"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "7.907"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "0.1628"
$ws.Range("E25").Value = "  +1.72%  "

$ws.Range("D26").Value = "163.85"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("D27").Value = "9.037"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("D28").Value = "18.57"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("E29").Value = "  +0.58%  "

$ws.Range("D30").Value = "4.411"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "1.289"
$ws.Range("E31").Value = "  -4.47%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "4.318"
$ws.Range("E32").Value = "  +4.95%  "

$ws.Range("D33").Value = "0.05254"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("D34").Value = "1.930"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("D35").Value = "0.7571"
$ws.Range("E35").Value = "  +3.85%  "

$ws.Range("D36").Value = "1.178"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").Value = "2.677"
$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").Value = "2.722"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").Value = "1.162.99"
$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D41").Value = "6.363"
$ws.Range("E41").Value = "  +4.10%  "

$ws.Range("D42").Value = "73.70"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("D44").Value = "104.87"
$ws.Range("E44").Value = "  +2.89%  "

$ws.Range("D45").Value = "0.9993"

$ws.Range("D46").Value = "2.031.33"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "1.800"
$ws.Range("E47").Value = "  +1.15%  "

$ws.Range("E48").Value = "  -1.59%  "

$ws.Range("D49").Value = "9.424"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.4315"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "7.050"
$ws.Range("E51").Value = "  +0.22%  "

# Restore normal (default) style on column D so no stray number-format styling remains
$dRange.Style = "Normal"
